$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-26 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-27 Friday", 2) | Out-Null
$d.Content.Find.Execute("935÷2=467, 1", $true, $false, $false, $false, $false, $true, 1, $false, "628÷9=69, 7", 2) | Out-Null
$d.Content.Find.Execute("682÷9=75, 7", $true, $false, $false, $false, $false, $true, 1, $false, "926÷4=231, 2", 2) | Out-Null
$d.Content.Find.Execute("393÷6=65, 3", $true, $false, $false, $false, $false, $true, 1, $false, "908÷8=113, 4", 2) | Out-Null
$d.Content.Find.Execute("394÷7=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "977÷7=139, 4", 2) | Out-Null
$d.Content.Find.Execute("855÷8=106, 7", $true, $false, $false, $false, $false, $true, 1, $false, "752÷6=125, 2", 2) | Out-Null
$d.Content.Find.Execute("880÷5=176, 0", $true, $false, $false, $false, $false, $true, 1, $false, "104÷8=13, 0", 2) | Out-Null
$d.Content.Find.Execute("492÷2=246, 0", $true, $false, $false, $false, $false, $true, 1, $false, "938÷5=187, 3", 2) | Out-Null
$d.Content.Find.Execute("233÷6=38, 5", $true, $false, $false, $false, $false, $true, 1, $false, "650÷2=325, 0", 2) | Out-Null
$d.Content.Find.Execute("641÷2=320, 1", $true, $false, $false, $false, $false, $true, 1, $false, "222÷4=55, 2", 2) | Out-Null
$d.Content.Find.Execute("419÷8=52, 3", $true, $false, $false, $false, $false, $true, 1, $false, "656÷9=72, 8", 2) | Out-Null
$d.Content.Find.Execute("328÷7=46, 6", $true, $false, $false, $false, $false, $true, 1, $false, "960÷4=240, 0", 2) | Out-Null
$d.Content.Find.Execute("584÷7=83, 3", $true, $false, $false, $false, $false, $true, 1, $false, "692÷9=76, 8", 2) | Out-Null
$d.Content.Find.Execute("809÷3=269, 2", $true, $false, $false, $false, $false, $true, 1, $false, "839÷9=93, 2", 2) | Out-Null
$d.Content.Find.Execute("952÷7=136, 0", $true, $false, $false, $false, $false, $true, 1, $false, "369÷9=41, 0", 2) | Out-Null
$d.Content.Find.Execute("502÷9=55, 7", $true, $false, $false, $false, $false, $true, 1, $false, "984÷8=123, 0", 2) | Out-Null
$d.Content.Find.Execute("743÷6=123, 5", $true, $false, $false, $false, $false, $true, 1, $false, "150÷4=37, 2", 2) | Out-Null
$d.Content.Find.Execute("707÷9=78, 5", $true, $false, $false, $false, $false, $true, 1, $false, "953÷3=317, 2", 2) | Out-Null
$d.Content.Find.Execute("228÷3=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "596÷2=298, 0", 2) | Out-Null
$d.Content.Find.Execute("233÷5=46, 3", $true, $false, $false, $false, $false, $true, 1, $false, "635÷7=90, 5", 2) | Out-Null
$d.Content.Find.Execute("491÷7=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "671÷7=95, 6", 2) | Out-Null
$d.Content.Find.Execute("557÷5=111, 2", $true, $false, $false, $false, $false, $true, 1, $false, "861÷2=430, 1", 2) | Out-Null
$d.Content.Find.Execute("936÷7=133, 5", $true, $false, $false, $false, $false, $true, 1, $false, "768÷9=85, 3", 2) | Out-Null
$d.Content.Find.Execute("551÷8=68, 7", $true, $false, $false, $false, $false, $true, 1, $false, "105÷6=17, 3", 2) | Out-Null
$d.Content.Find.Execute("425÷7=60, 5", $true, $false, $false, $false, $false, $true, 1, $false, "291÷9=32, 3", 2) | Out-Null
$d.Content.Find.Execute("500÷8=62, 4", $true, $false, $false, $false, $false, $true, 1, $false, "222÷8=27, 6", 2) | Out-Null
